# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (total) sheet,
#    populated with the per-fund holding breakdown for 2022-Q1.
# 2. Update the "总计" sheet with a new leading row summarising 2022-Q1
#    (holding count + market value), pushing the existing quarters down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Locate the existing "总计" sheet (currently the last tab) so the new
# sheet can be inserted immediately before it.
# ---------------------------------------------------------------------
$totalSheetOld = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($totalSheetOld)
$newSheet.Name = "2022-Q1"

# NOTE: inserting a sheet shifts tab positions, so any handle obtained
# *before* the insert (e.g. $totalSheetOld) can no longer be trusted to
# refer to the same sheet afterwards. Re-resolve "总计" now that it has
# settled into its final (last) position.
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy a style-2 cell (bold / bordered / centered header style) from the
# "总计" sheet onto the new sheet's header row and index column, so the
# new sheet matches the look of the other quarterly sheets without
# creating brand-new style entries.
$totalSheet.Range("B1").Copy($newSheet.Range("B1:H1"))
$totalSheet.Range("A2").Copy($newSheet.Range("A2:A7"))

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns D, E, F, G hold numeric-looking figures that are stored as text
# in the source data, so force text number-formatting before assigning.
$newSheet.Range("D2:G7").NumberFormat = "@"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "506006"
$newSheet.Range("C2").Value = "汇添富科创板2年定期开放混合"
$newSheet.Range("D2").Value = "26.37"
$newSheet.Range("E2").Value = "91.69"
$newSheet.Range("F2").Value = "3.51"
$newSheet.Range("G2").Value = "0.9256"
$newSheet.Range("H2").Value = 8

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "002707"
$newSheet.Range("C3").Value = "摩根士丹利华鑫科技领先灵活配置混合"
$newSheet.Range("D3").Value = "2.27"
$newSheet.Range("E3").Value = "93.05"
$newSheet.Range("F3").Value = "7.21"
$newSheet.Range("G3").Value = "0.1637"
$newSheet.Range("H3").Value = 1

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "501201"
$newSheet.Range("C4").Value = "红土创新科技创新 3 年封闭运作灵活配置混合"
$newSheet.Range("D4").Value = "3.99"
$newSheet.Range("E4").Value = "96.70"
$newSheet.Range("F4").Value = "2.99"
$newSheet.Range("G4").Value = "0.1193"
$newSheet.Range("H4").Value = 9

# Row 5
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "011603"
$newSheet.Range("C5").Value = "兴业高端制造混合A"
$newSheet.Range("D5").Value = "1.19"
$newSheet.Range("E5").Value = "76.40"
$newSheet.Range("F5").Value = "2.86"
$newSheet.Range("G5").Value = "0.0340"
$newSheet.Range("H5").Value = 6

# Row 6
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "168401"
$newSheet.Range("C6").Value = "红土创新转型精选灵活配置混合（LOF）"
$newSheet.Range("D6").Value = "0.78"
$newSheet.Range("E6").Value = "93.82"
$newSheet.Range("F6").Value = "3.24"
$newSheet.Range("G6").Value = "0.0253"
$newSheet.Range("H6").Value = 10

# Row 7
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "011604"
$newSheet.Range("C7").Value = "兴业高端制造混合C"
$newSheet.Range("D7").Value = "0.54"
$newSheet.Range("E7").Value = "76.40"
$newSheet.Range("F7").Value = "2.86"
$newSheet.Range("G7").Value = "0.0154"
$newSheet.Range("H7").Value = 6

# ---------------------------------------------------------------------
# 2. Update "总计" sheet: push existing quarters down one row and add
#    the new 2022-Q1 summary row at the top of the data.
# ---------------------------------------------------------------------

# Give the new last data row (A7) the same index style as the rest of
# column A before writing into it.
$totalSheet.Range("A6").Copy($totalSheet.Range("A7"))

$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2020-Q4"
$totalSheet.Range("C7").Value = 9
$totalSheet.Range("D7").Value = 5.28

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q1"
$totalSheet.Range("C6").Value = 28
$totalSheet.Range("D6").Value = 9.449999999999999

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 38
$totalSheet.Range("D5").Value = 25.92

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 31
$totalSheet.Range("D4").Value = 26.6

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 38
$totalSheet.Range("D3").Value = 31.08

# New first data row (A2 already carries the index-cell style from the
# original data, so just overwrite its value).
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 1.28
